$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $r1, $r2, $cols) {
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}

# Columns that hold the taxon-observation specific data (not the shared
# location/date/observer columns) for rows 16/17 and 26/27.
$cols1617 = @("A","B","E","F","G","H","Q","R","AJ","AK","AO")
Swap-Rows $ws 16 17 $cols1617

$cols2627 = @("A","B","E","F","G","H","Q","R")
Swap-Rows $ws 26 27 $cols2627
